$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two more rows (4 and 5) that duplicate the existing rows 2 and 3,
# matching the new match entries added to the scraped data.
# Force the numeric-looking columns (G:K) to be stored as text, same as
# the rest of the sheet (numberStoredAsText), then restore the original
# (default) cell style so no new number-format style is left applied.
$ws.Range("G4:K5").NumberFormat = "@"

$ws.Range("A4").Value = " Sharjah"
$ws.Range("B4").Value = " October 23 2020"
$ws.Range("C4").Value = "Mumbai won by 10 wickets (with 46 balls remaining)"
$ws.Range("D4").Value = "Chennai Super Kings"
$ws.Range("E4").Value = "Mumbai Indians"
$ws.Range("F4").Value = "Shardul Thakur "
$ws.Range("G4").Value = "11"
$ws.Range("H4").Value = "20"
$ws.Range("I4").Value = "0"
$ws.Range("J4").Value = "0"
$ws.Range("K4").Value = "55.00"

$ws.Range("A5").Value = " Dubai (DSC)"
$ws.Range("B5").Value = " October 10 2020"
$ws.Range("C5").Value = "RCB won by 37 runs"
$ws.Range("D5").Value = "Chennai Super Kings"
$ws.Range("E5").Value = "Royal Challengers Bangalore"
$ws.Range("F5").Value = "Shardul Thakur "
$ws.Range("G5").Value = "1"
$ws.Range("H5").Value = "1"
$ws.Range("I5").Value = "0"
$ws.Range("J5").Value = "0"
$ws.Range("K5").Value = "100.00"

# The NumberFormat="@" step above stamped a "Text" style onto G4:K5; put
# the plain default style (matching the rest of the sheet) back on them.
$ws.Range("G4:K5").Style = $ws.Range("G2:K3").Style
